$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# --- Update the "Pthread interval 2000-5000" data table (columns B:D, rows 2-11) ---
# Column C ("Rect-2") and column D ("Trap.") values were re-measured; column B ("Rect-1")
# is unchanged except for the last row.

$ws.Range("C2").Value2 = 0.00100099999999999994
$ws.Range("D2").Value2 = 0.00200000000000000004

$ws.Range("C3").Value2 = 0.00300599999999999999
$ws.Range("D3").Value2 = 0.00200000000000000004

$ws.Range("C4").Value2 = 0.0029970000000000001
$ws.Range("D4").Value2 = 0.00399899999999999995

$ws.Range("C5").Value2 = 0.00399699999999999969
$ws.Range("D5").Value2 = 0.00499699999999999971

$ws.Range("C6").Value2 = 0.00400100000000000022
$ws.Range("D6").Value2 = 0.00400200000000000035

$ws.Range("C7").Value2 = 0.00499899999999999997
$ws.Range("D7").Value2 = 0.00498999999999999964

$ws.Range("C8").Value2 = 0.0059890000000000004
$ws.Range("D8").Value2 = 0.0050000000000000001

$ws.Range("C9").Value2 = 0.00500100000000000024
$ws.Range("D9").Value2 = 0.00400000000000000008

$ws.Range("C10").Value2 = 0.00599899999999999999
$ws.Range("D10").Value2 = 0.0050000000000000001

$ws.Range("B11").Value2 = 0.00500200000000000037
$ws.Range("C11").Value2 = 0.0050000000000000001
$ws.Range("D11").Value2 = 0.00600200000000000039

# --- Update the active selection left behind on the sheet ---
$ws.Range("H13").Select()
